# Update "planeamento de tarefas.xlsx" - Folha1 (sheet1)
# - Rename the three minigame tasks to include their specific names
# - Mark a few tasks as started ("iniciado") by Filipe (column H = Filipe, column I = Estado)
# - Update the selected cell in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename minigame task labels (appear twice: once under "navegação 3D" group,
# once under "inteligência artificial" group)
$ws.Range("B12").Value = "Minijogo 1 - jogo do galo"
$ws.Range("B13").Value = "Minijogo 2 - enforcado"
$ws.Range("B14").Value = "Minijogo 3 - labirinto"

$ws.Range("B19").Value = "Minijogo 1 - jogo do galo"
$ws.Range("B20").Value = "Minijogo 2 - enforcado"
$ws.Range("B21").Value = "Minijogo 3 - labirinto"

# Mark progress: Filipe ("x") and Estado ("iniciado") for a few tasks
$ws.Range("H5").Value = "x"
$ws.Range("I5").Value = "iniciado"

$ws.Range("H6").Value = "x"
$ws.Range("I6").Value = "iniciado"

$ws.Range("H8").Value = "x"
$ws.Range("I8").Value = "iniciado"

$ws.Range("I12").Value = "iniciado"

# Update selection / active cell shown when the sheet is reopened
$ws.Range("H12").Select()
